$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 5 more rows (rows 7-11) that duplicate the existing data rows (2-6).
$policyNumber = 573
$sourceUrl = "https://al-policies.exploremyplan.com/portal/web/medical-policies/-/mp-573"
$destUrl = "https://stage-us-mypolicies.itilitihealth.us/policy/938125692074/573?lob=BCBS+AL"

for ($r = 7; $r -le 11; $r++) {
    $ws.Cells.Item($r, 1).Value = $policyNumber
    $ws.Cells.Item($r, 2).Value = $sourceUrl
    $ws.Cells.Item($r, 3).Value = $destUrl
}
